# Fix misaligned course-syllabus rows (LOQ4222) and insert missing content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout: split column A (label) off from column B (PT content) ---
# Column B must end up the same width as column C (60.7109375 chars); nudge it
# so the engine re-splits the <col> run that used to cover A:B together.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Row 10
$ws.Range("B10").Value = 'Apresentar noções de Matemática Financeira, Gestão Financeiras e Engenharia Econômica aos alunos, capacitando-os para construir e analisar fluxos de caixa de projetos e empreendimentos, discutir os principais aspectos da gestão financeira das empresas industriais, comerciais e de serviços e analisar e propor estratégias de gestão financeira relacionadas às estratégias de mercado e de produção.'
$ws.Range("C10").Value = 'Apresentar noções de Matemática Financeira, Gestão Financeiras e Engenharia Econômica aos alunos, capacitando-os para construir e analisar fluxos de caixa de projetos e empreendimentos, discutir os principais aspectos da gestão financeira das empresas industriais, comerciais e de serviços e analisar e propor estratégias de gestão financeira relacionadas às estratégias de mercado e de produção.'

# Row 13
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '1. Engenharia Econômica. 2 – Finanças'
$ws.Range("C14").Value = '1. Engenharia Econômica. 2 – Finanças'

# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = '1. Economic Engineering. 2 – Finance'
$ws.Range("C15").Value = '1. Economic Engineering. 2 – Finance'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1. Engenharia Econômica: Variável tempo: juros simples, juros compostos; Métodos de amortização; Equivalência de métodos; Métodos de Decisão; Renovação e substituição de equipamentos; Depreciação; Análise de Projetos, Riscos em projetos; Estimativa do custo de capital próprio (CAPM) e WACC.2. Finanças: O ciclo da produção e o ciclo do capital; Análise de Índices; Fontes de Financiamento, Alavancagem; Capital de Giro; Custo de Capital; Ações, Política de Dividendos; Financiamento de Longo Prazo, Corporate Finance/Project Finance; EVA e MVA.'
$ws.Range("C16").Value = '1. Engenharia Econômica: Variável tempo: juros simples, juros compostos; Métodos de amortização; Equivalência de métodos; Métodos de Decisão; Renovação e substituição de equipamentos; Depreciação; Análise de Projetos, Riscos em projetos; Estimativa do custo de capital próprio (CAPM) e WACC.2. Finanças: O ciclo da produção e o ciclo do capital; Análise de Índices; Fontes de Financiamento, Alavancagem; Capital de Giro; Custo de Capital; Ações, Política de Dividendos; Financiamento de Longo Prazo, Corporate Finance/Project Finance; EVA e MVA.'

# Row 17
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1. Economic Engineering: Time variable: simple interest, compound interest; Amortization methods; Equivalence of methods; Decision Methods; Renovation and replacement of equipment; Depreciation; Project Analysis, Project Risks; Estimated cost of equity (CAPM) and WACC.2. Finance: The production cycle and the capital cycle; Index Analysis; Financing Sources, Leverage; Working capital; Capital cost; Shares, Dividend Policy; Long Term Financing, Corporate Finance/Project Finance; EVA and MVA.'
$ws.Range("C17").Value = '1. Economic Engineering: Time variable: simple interest, compound interest; Amortization methods; Equivalence of methods; Decision Methods; Renovation and replacement of equipment; Depreciation; Project Analysis, Project Risks; Estimated cost of equity (CAPM) and WACC.2. Finance: The production cycle and the capital cycle; Index Analysis; Financing Sources, Leverage; Working capital; Capital cost; Shares, Dividend Policy; Long Term Financing, Corporate Finance/Project Finance; EVA and MVA.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("A18").Value = 'Avaliação:'
$ws.Rows.Item(18).AutoFit()

# Row 19
$ws.Range("A19").Value = 'Método:'

# Row 20
$ws.Range("A20").Value = 'Critério:'

# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRAÇÃO FINANCEIRA. São Paulo: Atlas, 2014MARIANO, F.; MENESES, A. Curso De Administração Financeira. São Paulo: Método, 2012.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econômica. São Paulo: LTC,2000.KOPITTKE, B. H.; CASAROTTO FILHO, N. ANÁLISE DE INVESTIMENTOS: Matemática Financeira, Engenharia Econômica, Estratégia Empresarial. 11 ed. São Paulo: Atlas, 2010.HOJI, M.; LUZ, A. E. Gestão Financeira Econômica: Didática, Objetiva e Prática. São Paulo: Atlas, 2019.GOMES, J. M. Elaboração e Análise De Viabilidade Econômica De Projetos. São Paulo: Atlas, 2013.OLIVO, R. L. F. Análise de Investimentos. Campinas: Alínea, 2011.ALMEIDA, J. T. S. Matemática Financeira. Rio de Janeiro: LTC, 2016.EHRLICH, Pierre Jacques. Engenharia Econômica. São Paulo: Editora Atlas, 2005.HIRSCHFELD, Henrique. Engenharia econômica e análise de custos. 7. ed. São Paulo: Atlas, 2007.MOTTA, Regis da Rocha; CALÔBA, Guilherme Marques. Análise de Investimentos. São Paulo: Atlas 2002.SANVICENTE, A. Z. Administração Financeira. São Paulo: Editora Atlas, 2007.VAN HORNE, J. C. Política e Administração Financeira. Rio de Janeiro: Livros Técnicos e Científicos, 1974.WESTON, J. F.; BRIGHAM, E. F. Administração Financeira de Empresas. São Paulo: Editora Interamericana, 2000.'
$ws.Range("C22").Value = 'GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRAÇÃO FINANCEIRA. São Paulo: Atlas, 2014MARIANO, F.; MENESES, A. Curso De Administração Financeira. São Paulo: Método, 2012.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econômica. São Paulo: LTC,2000.KOPITTKE, B. H.; CASAROTTO FILHO, N. ANÁLISE DE INVESTIMENTOS: Matemática Financeira, Engenharia Econômica, Estratégia Empresarial. 11 ed. São Paulo: Atlas, 2010.HOJI, M.; LUZ, A. E. Gestão Financeira Econômica: Didática, Objetiva e Prática. São Paulo: Atlas, 2019.GOMES, J. M. Elaboração e Análise De Viabilidade Econômica De Projetos. São Paulo: Atlas, 2013.OLIVO, R. L. F. Análise de Investimentos. Campinas: Alínea, 2011.ALMEIDA, J. T. S. Matemática Financeira. Rio de Janeiro: LTC, 2016.EHRLICH, Pierre Jacques. Engenharia Econômica. São Paulo: Editora Atlas, 2005.HIRSCHFELD, Henrique. Engenharia econômica e análise de custos. 7. ed. São Paulo: Atlas, 2007.MOTTA, Regis da Rocha; CALÔBA, Guilherme Marques. Análise de Investimentos. São Paulo: Atlas 2002.SANVICENTE, A. Z. Administração Financeira. São Paulo: Editora Atlas, 2007.VAN HORNE, J. C. Política e Administração Financeira. Rio de Janeiro: Livros Técnicos e Científicos, 1974.WESTON, J. F.; BRIGHAM, E. F. Administração Financeira de Empresas. São Paulo: Editora Interamericana, 2000.'
$ws.Rows.Item(22).RowHeight = 120

# Row 23
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Range("A23").Value = 'Requisitos:'
$ws.Rows.Item(23).AutoFit()

# Row 24
$ws.Range("B24").Value = 'LOQ4213 -  Contabilidade e Custos  (Requisito fraco)`n'
$ws.Range("C24").Value = 'LOQ4213 -  Contabilidade e Custos  (Requisito fraco)`n'
$ws.Rows.Item(24).RowHeight = 30
